# "after incorporating data from more recent studies"
#
# Changes applied to Sheet 1:
#  1. Rename the author cell "Sun et al. (2023)" (row 8, column A) to
#     "Sun (2023)" - its correlation values are unchanged.
#  2. Insert two new rows of data for newly-incorporated studies right
#     before the existing "Baez et al. (2017)" row (which was row 9 and
#     is pushed down to row 11):
#       Brog et al. (2022) -> Treatment 0.137092597312419, Control 0.994257446012105
#       Iyer et al. (2023) -> Treatment 0.0929704290319567, Control 0.991422978792847

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename "Sun et al. (2023)" -> "Sun (2023)" (row 8, values unaffected).
$ws.Cells.Item(8, 1).Value = "Sun (2023)"

# 2. Insert two blank rows at row 9, pushing "Baez et al. (2017)" (and
#    anything after it) down by two rows.
$ws.Rows.Item(9).Resize(2).EntireRow.Insert()

# New row 9: Brog et al. (2022)
$ws.Cells.Item(9, 1).Value = "Brog et al. (2022)"
$ws.Cells.Item(9, 2).Value = 0.137092597312419
$ws.Cells.Item(9, 3).Value = 0.994257446012105

# New row 10: Iyer et al. (2023)
$ws.Cells.Item(10, 1).Value = "Iyer et al. (2023)"
$ws.Cells.Item(10, 2).Value = 0.0929704290319567
$ws.Cells.Item(10, 3).Value = 0.991422978792847
